$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3911
$ws1.Range("F4").Value = 2309
$ws1.Range("F7").Value = 26
$ws1.Range("F8").Value = 188
$ws1.Range("F11").Value = 1445
$ws1.Range("F13").Value = 2599

# Sheet "全部类型" (All Types) - same metric, mirrored rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3911
$ws4.Range("F4").Value = 2309
$ws4.Range("F7").Value = 26
$ws4.Range("F9").Value = 188
$ws4.Range("F14").Value = 1445
$ws4.Range("F16").Value = 2599
